$d = $word.ActiveDocument

# 1. Fix typo: "BlosHGome" -> "BlosHome"
$d.Content.Find.Execute("BlosHGome", $false, $false, $false, $false, $false,
                         $true, 1, $false, "BlosHome", 2)

# 2. "...avec bonheur. Mes deux blogs" -> "...avec bonheur. Mes deux"
#    (drop the trailing " blogs" that belonged to the first run only -
#    the next run already supplies its own " "/"blogs " text, untouched)
$r2 = $d.Content
$r2.Find.Execute("avec bonheur. Mes deux blogs")
$suffix = " blogs"
$del = $d.Range($r2.End - $suffix.Length, $r2.End)
$del.Text = ""

# 3. "and " -> "et" + " " (split into two runs with identical formatting)
$r3 = $d.Content
$r3.Find.Execute("and ")
$r3.Text = "et"
$r3.Collapse(0)
$r3.InsertAfter(" ")
